$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = "13 (2.6) "
$ws.Range("B3").Value = "4.7 (1.3) "
$ws.Range("B4").Value = "42 (1.7) "
$ws.Range("B5").Value = "36 (4.0) "
$ws.Range("B6").Value = "21 (5.1) "
$ws.Range("B7").Value = "68 (16) "
$ws.Range("F7").Select()
